$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-34 (Generation 0-32): Fitness column C changes from 7310 to 7295
$ws.Range("C2:C34").Value = 7295

# Rows 35-252 (Generation 33-250): Fitness column C changes from 7310 to 7293
$ws.Range("C35:C252").Value = 7293
